$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "(312049950, Molham  Peretz: 3,6)"
$ws.Range("B1").Value = "(308073899, Anan  Kirshenbaum: 7,6)"
$ws.Range("C1").Value = "(318869187, Soaad  Leibovich: -4,2)"
$ws.Range("D1").Value = "(205898513, Asaf  Braymok: -8,-4)"
$ws.Range("E1").Value = "(318428158, Tal  Asulin: 3,4)"
$ws.Range("F1").Value = "(316028364, Sami  Castro: 8,9)"
$ws.Range("G1").Value = "(318294931, Shalev  Afanasenko: 2,-7)"

$ws.Range("A3").Value = "cost: 303.1496768052041"
$ws.Range("A4").Value = "time: 55.62993536104082"
